# Auto-generated: appends new sensor log rows to PIR, Humidity, Temperature, Proximity, and Camera sheets
function Add-LogRow($ws, $r, $a, $b, $c, $d, $e, $f) {
    # Force the whole row to Text format first so Excel doesn't
    # auto-convert date-looking, percent-looking, or numeric-looking
    # strings (e.g. "2026-01-28", "88.4%") into dates/numbers - the
    # source log stores every cell as plain text.
    $ws.Range("A$r`:F$r").NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
}

$wb = $excel.ActiveWorkbook

# --- PIR sheet ---
$ws = $wb.Worksheets.Item("PIR")
Add-LogRow $ws 294 "2026-01-28" "15:16:36" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 295 "2026-01-28" "15:16:38" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 296 "2026-01-28" "15:16:43" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 297 "2026-01-28" "15:16:48" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 298 "2026-01-28" "15:16:54" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 299 "2026-01-28" "15:16:59" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 300 "2026-01-28" "15:17:04" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 301 "2026-01-28" "15:17:10" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 302 "2026-01-28" "15:17:14" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 303 "2026-01-28" "15:17:19" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 304 "2026-01-28" "15:17:24" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 305 "2026-01-28" "15:17:29" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 306 "2026-01-28" "15:17:34" "15:00" "Bathroom" "No Motion" "Inactive"

# --- Humidity sheet ---
$ws = $wb.Worksheets.Item("Humidity")
Add-LogRow $ws 280 "2026-01-28" "15:16:34" "15:00" "Bathroom" "88.4%" "Active"
Add-LogRow $ws 281 "2026-01-28" "15:16:40" "15:00" "Bathroom" "87.4%" "Active"
Add-LogRow $ws 282 "2026-01-28" "15:16:52" "15:00" "Bathroom" "88.3%" "Active"
Add-LogRow $ws 283 "2026-01-28" "15:16:56" "15:00" "Bathroom" "88.3%" "Active"
Add-LogRow $ws 284 "2026-01-28" "15:17:04" "15:00" "Bathroom" "88.3%" "Active"
Add-LogRow $ws 285 "2026-01-28" "15:17:08" "15:00" "Bathroom" "87.4%" "Active"
Add-LogRow $ws 286 "2026-01-28" "15:17:12" "15:00" "Bathroom" "88.3%" "Active"
Add-LogRow $ws 287 "2026-01-28" "15:17:16" "15:00" "Bathroom" "88.3%" "Active"
Add-LogRow $ws 288 "2026-01-28" "15:17:20" "15:00" "Bathroom" "87.4%" "Active"
Add-LogRow $ws 289 "2026-01-28" "15:17:25" "15:00" "Bathroom" "88.3%" "Active"
Add-LogRow $ws 290 "2026-01-28" "15:17:32" "15:00" "Bathroom" "87.4%" "Active"

# --- Temperature sheet ---
$ws = $wb.Worksheets.Item("Temperature")
Add-LogRow $ws 280 "2026-01-28" "15:16:35" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 281 "2026-01-28" "15:16:41" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 282 "2026-01-28" "15:16:53" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 283 "2026-01-28" "15:16:57" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 284 "2026-01-28" "15:17:05" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 285 "2026-01-28" "15:17:09" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 286 "2026-01-28" "15:17:13" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 287 "2026-01-28" "15:17:17" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 288 "2026-01-28" "15:17:21" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 289 "2026-01-28" "15:17:25" "15:00" "Bathroom" "22.9C" "Active"
Add-LogRow $ws 290 "2026-01-28" "15:17:33" "15:00" "Bathroom" "22.9C" "Active"

# --- Proximity sheet ---
$ws = $wb.Worksheets.Item("Proximity")
Add-LogRow $ws 41 "2026-01-28" "15:17:30" "15:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"

# --- Camera sheet ---
$ws = $wb.Worksheets.Item("Camera")
Add-LogRow $ws 19 "2026-01-28" "15:17:31" "15:00" "Living Room Main Door" "Image Captured" "Active"
